$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
#    Replace the full howstat URL with just the numeric match code, and drop
#    the handful of cells in column B that only held an empty inline string.
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$battingUsed = $batting.UsedRange
$battingLastRow = $battingUsed.Row + $battingUsed.Rows.Count - 1

$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

# Format column D as text first so the bare numeric codes are stored as
# strings (matching the source data) rather than being coerced to numbers.
$batting.Range("D2:D$battingLastRow").NumberFormat = "@"

for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $link = [string]$cell.Value2
    if ($link -match 'MatchCode=(\d+)') {
        $cell.Value = $matches[1]
    }
}

# Rows whose INNING_NUMBER cell is an empty placeholder ("") get cleared out
# entirely so the <c> element disappears, as in the source edit.
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, 2)
    $v = $cell.Value2
    if ($v -eq "" -or $v -eq $null) {
        $cell.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 2) Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowlingUsed = $bowling.UsedRange
$bowlingLastRow = $bowlingUsed.Row + $bowlingUsed.Rows.Count - 1

$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowling.Range("B2:B$bowlingLastRow").NumberFormat = "@"

for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $link = [string]$cell.Value2
    if ($link -match 'MatchCode=(\d+)') {
        $cell.Value = $matches[1]
    }
}

# ---------------------------------------------------------------------------
# 3) Add new "Player Info" sheet before "ODI Batting"
# ---------------------------------------------------------------------------
$battingRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingRef, $null)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $playerInfoHeaders[$c - 1]
}
$batting.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfoRow = @("4266", "Mohammad Mithun Ali", "Right Handed", "Does Not Bowl | Unknown")
for ($c = 1; $c -le $playerInfoRow.Length; $c++) {
    $playerInfo.Cells.Item(2, $c).Value = $playerInfoRow[$c - 1]
}

# ---------------------------------------------------------------------------
# 4) Add new "ODI Batting Extra" sheet after "ODI Bowling"
# ---------------------------------------------------------------------------
$bowlingRef = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingRef)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $extra.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
$batting.Range("A1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

$extraRows = @(
    @("4251", 5, "7", "1", "25.22%", "NO"),
    @("4286", 5, $null, $null, $null, "NO"),
    @("4293", 5, "2", "2", "17.34%", "NO"),
    @("4296", 5, "1", "1", "7.98%", "NO"),
    @("4307", $null, $null, $null, $null, "NO"),
    @("4311", $null, $null, $null, $null, "NO"),
    @("4314", 5, "0", "0", $null, "NO"),
    @("4356", $null, $null, $null, $null, "NO"),
    @("4357", $null, $null, $null, $null, "NO"),
    @("4358", $null, $null, $null, $null, "NO"),
    @("4416", 6, "5", "1", "15.58%", "NO"),
    @("4418", $null, $null, $null, $null, "NO"),
    @("4420", 5, $null, $null, $null, "NO"),
    @("4452", 5, "1", "0", "6.87%", "NO"),
    @("4453", 5, "6", "2", "26.94%", "NO"),
    @("4455", 4, "0", "0", "3.90%", "NO"),
    @("4463", 5, "0", "0", $null, "NO"),
    @("4477", 4, "4", "0", "6.88%", "NO"),
    @("4479", 4, "0", "0", "0.83%", "NO"),
    @("4481", 4, "1", "0", "9.93%", "NO")
)

$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:C21").NumberFormat = "@"
$extra.Range("D2:D21").NumberFormat = "@"
$extra.Range("E2:E21").NumberFormat = "@"

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne $null) {
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne $null) {
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}

Write-Output "done"
